$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f6ae490>),
                (''model'',
                 LogisticRegression(l1_ratio=0.7, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Range("B2").Value = 0.7285714285714284
$ws.Range("C2").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f537b20>, ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.7, ''model__class_weight'': None}'
$ws.Range("D2").Value = 0.7030254444431433
$ws.Range("E2").Value = 0.6415745712620712
$ws.Range("F2").Value = 0.8
$ws.Range("G2").Value = 0.6766131624076245
$ws.Range("H2").Value = 0.6188425925925926
$ws.Range("I2").Value = 0.6666666666666666
$ws.Range("J2").Value = 0.7644060283687942
$ws.Range("K2").Value = 0.7086805555555555
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f537040>),
                (''model'',
                 LogisticRegression(l1_ratio=0.95, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Range("B3").Value = 0.7154761904761904
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f537970>, ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.95, ''model__class_weight'': None}'
$ws.Range("D3").Value = 0.6974039380237601
$ws.Range("E3").Value = 0.5698129532504531
$ws.Range("F3").Value = 0.8
$ws.Range("G3").Value = 0.6734361846925557
$ws.Range("H3").Value = 0.6046288029100528
$ws.Range("I3").Value = 0.6666666666666666
$ws.Range("J3").Value = 0.7523492907801419
$ws.Range("K3").Value = 0.5895833333333332
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f537940>),
                (''model'',
                 LogisticRegression(class_weight=''balanced'', l1_ratio=0.01,
                                    max_iter=1000, penalty=''elasticnet'',
                                    random_state=42, solver=''saga''))])'
$ws.Range("B4").Value = 0.678124098124098
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f53bc10>, ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.01, ''model__class_weight'': ''balanced''}'
$ws.Range("D4").Value = 0.6759563493888675
$ws.Range("E4").Value = 0.5868948875198875
$ws.Range("F4").Value = 0.5714285714285714
$ws.Range("G4").Value = 0.6779862284556121
$ws.Range("H4").Value = 0.6209151785714285
$ws.Range("I4").Value = 0.8888888888888888
$ws.Range("J4").Value = 0.6875
$ws.Range("K4").Value = 0.5979166666666667
$ws.Range("L4").Value = 0.4210526315789473
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[0 1 1 0 1 1 0 1 1 0 0 0 0 0 1 0 0 1 0 0 0 0 1 0]'
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f53bb80>),
                (''model'',
                 LogisticRegression(l1_ratio=0.5, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Range("B5").Value = 0.7499999999999999
$ws.Range("C5").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f53b430>, ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.5, ''model__class_weight'': None}'
$ws.Range("D5").Value = 0.7143924393270283
$ws.Range("E5").Value = 0.658300992988493
$ws.Range("F5").Value = 0.7368421052631579
$ws.Range("G5").Value = 0.6772385259720881
$ws.Range("H5").Value = 0.6225231481481481
$ws.Range("I5").Value = 0.5833333333333334
$ws.Range("J5").Value = 0.7948979591836736
$ws.Range("K5").Value = 0.7531249999999999
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f53b310>),
                (''model'',
                 LogisticRegression(l1_ratio=0.5, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Range("B6").Value = 0.7416666666666666
$ws.Range("C6").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f775e20>, ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.5, ''model__class_weight'': None}'
$ws.Range("D6").Value = 0.7697580797446871
$ws.Range("E6").Value = 0.6596079846079845
$ws.Range("F6").Value = 0.6285714285714286
$ws.Range("G6").Value = 0.732290848301124
$ws.Range("H6").Value = 0.6331068121693122
$ws.Range("I6").Value = 0.4583333333333333
$ws.Range("J6").Value = 0.8368589743589743
$ws.Range("K6").Value = 0.7288194444444445
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'